$wb = $excel.ActiveWorkbook

# Update both the "展览" and "全部类型" sheets (they hold duplicate data)
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1739
    $ws.Range("F3").Value = 7977
    $ws.Range("F5").Value = 277
}
